# Apply the table-style change to the three tables that used the custom
# "Table_0" style, switching them to the built-in table style
# {848B3430-8FF2-4F73-AEBE-68E9FE146709}.
$p = $ppt.ActivePresentation

$newTableStyleId = "{848B3430-8FF2-4F73-AEBE-68E9FE146709}"

foreach ($slideIndex in 14, 15, 16) {
    $slide = $p.Slides.Item($slideIndex)
    for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
        $shape = $slide.Shapes.Item($i)
        if ($shape.HasTable) {
            $shape.Table.ApplyStyle($newTableStyleId)
        }
    }
}

# Swap the presentation's theme colour scheme from the custom
# "Red Violet" / Integral palette over to the standard Office Theme
# palette (the two theme parts effectively traded their colour schemes).
$officeColors = @(
    "000000",
    "FFFFFF",
    "44546A",
    "E7E6E6",
    "5B9BD5",
    "ED7D31",
    "A5A5A5",
    "FFC000",
    "4472C4",
    "70AD47",
    "0563C1",
    "954F72"
)

$slide1 = $p.Slides.Item(1)
$colorScheme = $slide1.ThemeColorScheme

for ($i = 1; $i -le $officeColors.Length; $i++) {
    $hex = $officeColors[$i - 1]
    $r = [Convert]::ToInt32($hex.Substring(0, 2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2, 2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4, 2), 16)
    $colorScheme.Item($i).RGB = $r + ($g * 256) + ($b * 65536)
}
